# Improved loghandling of Write To Dynamics CRM target (logs now filterable)
#
# This script reproduces, via the Excel COM object model, the edits that were
# captured in the OOXML diff:
#   - Two new contact rows (7 and 8) are appended to the "Contacts" sheet,
#     with two new shared strings ("James again" / "Peter again").
#   - The "Contacts" sheet becomes the active/selected sheet (selection moves
#     to B6), while "InvoiceContacts" is no longer the selected sheet.

$wb = $excel.ActiveWorkbook

$contacts = $wb.Worksheets.Item("Contacts")

# --- New row 7: ContactId 4 / James again / Bond / A-0003 / U1 / '10.02.1989 (text) ---
$contacts.Range("A7").Value = 4
$contacts.Range("B7").Value = "James again"
$contacts.Range("C7").Value = "Bond"
$contacts.Range("D7").Value = "A-0003"
$contacts.Range("E7").Value = "U1"

# Reuse the text/quote-prefixed date format already used by F5 ("10.02.1989")
$contacts.Range("F5").Copy()
$contacts.Range("F7").PasteSpecial(-4122)
$contacts.Range("F7").Value = "'10.02.1989"

# --- New row 8: ContactId 2 / Peter again / Chan / A-0004 / U2 / 01.12.1989 (date) ---
$contacts.Range("A8").Value = 2
$contacts.Range("B8").Value = "Peter again"
$contacts.Range("C8").Value = "Chan"
$contacts.Range("D8").Value = "A-0004"
$contacts.Range("E8").Value = "U2"

# Reuse the real-date format already used by F3 (12/01/1989)
$contacts.Range("F3").Copy()
$contacts.Range("F8").PasteSpecial(-4122)
$contacts.Range("F8").Value = "12/01/1989"

# --- Sheet/selection activation ---
# "Contacts" becomes the selected tab with B6 selected ...
$contacts.Activate()
$contacts.Range("B6").Select()

# ... and "InvoiceContacts" (previously tabSelected) keeps its own selection
# but is no longer the active tab.
$invoiceContacts = $wb.Worksheets.Item("InvoiceContacts")
$invoiceContacts.Range("B2").Select()
$contacts.Activate()
